$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 300 (Excel-style: row 300's EntireRow.Insert
# shifts row 300..380 down to 301..381, matching the xlShiftDown default).
$ws.Rows.Item(300).Insert()

# Populate the newly inserted row 300 with the new data record.
$ws.Cells.Item(300, 1).Value  = 10
$ws.Cells.Item(300, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(300, 3).Value  = "La Araucanía"
$ws.Cells.Item(300, 4).Value  = 44841
$ws.Cells.Item(300, 5).Value  = 9
$ws.Cells.Item(300, 6).Value  = 100112017
$ws.Cells.Item(300, 7).Value  = "Apio"
$ws.Cells.Item(300, 8).Value  = "Americana (o)"
$ws.Cells.Item(300, 9).Value  = "Primera"
$ws.Cells.Item(300, 10).Value = 50
$ws.Cells.Item(300, 11).Value = 9000
$ws.Cells.Item(300, 12).Value = 9000
$ws.Cells.Item(300, 13).Value = 9000
$ws.Cells.Item(300, 14).Value = "`$/docena de matas"
$ws.Cells.Item(300, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(300, 16).Value = 1500
$ws.Cells.Item(300, 17).Value = 6
$ws.Cells.Item(300, 18).Value = "Hortaliza"

# Match the date-number formatting used by the rest of column D.
$ws.Cells.Item(300, 4).NumberFormat = $ws.Cells.Item(301, 4).NumberFormat
